# Adds two new worksheets (AddACHMRF, AddModifyDeleteACH) with ACH test
# data, mirroring the existing CreateProfile / AddModifyDeleteCC sheets,
# and tweaks the previously-active sheet's view now that it is no longer
# the selected tab.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) AddModifyDeleteCC is no longer the active/selected tab - its saved
#    view moves back to a neutral selection.
# ---------------------------------------------------------------------
$addModifyDeleteCC = $wb.Worksheets.Item("AddModifyDeleteCC")
$addModifyDeleteCC.Activate()
$addModifyDeleteCC.Application.ActiveWindow.ScrollColumn = 34
$addModifyDeleteCC.Range("AA1").Select() | Out-Null
$addModifyDeleteCC.Range("A1:XFD1048576").Select() | Out-Null

# ---------------------------------------------------------------------
# 2) New sheet "AddACHMRF" - a CreateProfile-style MRF sheet seeded with
#    the QA/Demo app ids + the default address block.
# ---------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$achMrf = $wb.Worksheets.Add($null, $lastSheet)
$achMrf.Name = "AddACHMRF"

$achMrfHeaders = @("Result","Date","Notes","Execute","AppIDQA","AppNameQA","AppIDDemo","AppNameDemo",
  "ProfileName","CompanyName","Title","FirstName","MiddleName","LastName","Suffix","AL1","AL2","Country",
  "ZIP","ZipExt","City","State","Email","HomePhone1","HomePhone2","HomePhone3","MobilePhone1","MobilePhone2",
  "MobilePhone3","WorkPhone1","WorkPhone2","WorkPhone3","WorkPhoneExt","Comments")
for ($i = 0; $i -lt $achMrfHeaders.Length; $i++) {
    $achMrf.Cells.Item(1, $i + 1).Value = $achMrfHeaders[$i]
}

$achMrf.Range("D2").Value = "Y"
$achMrf.Range("E2").Value = "742"
$achMrf.Range("F2").Value = "a_Access AutoNoCFtp"
$achMrf.Range("G2").Value = "4249"
$achMrf.Range("H2").Value = "a_Access AutoNoCFtpDemo"
$achMrf.Range("P2").Value = "15 Toledo Road"
$achMrf.Range("R2").Value = "United States"
$achMrf.Range("S2").Value = "22201"

$achMrf.Range("U1").Select() | Out-Null
$achMrf.Range("F8").Select() | Out-Null

# ---------------------------------------------------------------------
# 3) New sheet "AddModifyDeleteACH" - copy of the AddModifyDeleteCC
#    layout (cols A:AH) plus the ACH-specific columns AI:AM, populated
#    with three rows (Personal Checking / Personal Savings / Business
#    Checking) of test data.
# ---------------------------------------------------------------------
$ach = $wb.Worksheets.Add($null, $achMrf)
$ach.Name = "AddModifyDeleteACH"

$achHeaders = @("Result","Date","Notes","Execute","AppIDQA","AppNameQA","AppIDDemo","AppNameDemo",
  "ProfileName","CompanyName","Title","FirstName","MiddleName","LastName","Suffix","AL1","AL2","Country",
  "ZIP","ZipExt","City","State","Email","HomePhone1","HomePhone2","HomePhone3","MobilePhone1","MobilePhone2",
  "MobilePhone3","WorkPhone1","WorkPhone2","WorkPhone3","WorkPhoneExt","Comments",
  "PaymentType","RTN","NicknameACH","ACNumber","Required Fields ACH")
for ($i = 0; $i -lt $achHeaders.Length; $i++) {
    $ach.Cells.Item(1, $i + 1).Value = $achHeaders[$i]
}

$achRows = @(
    @{ Type = "Personal Checking"; Nick = "Ross PC"; Acct = "25872222"; NickMod = "Ross PC Mod" },
    @{ Type = "Personal Savings"; Nick = "Ross PS"; Acct = "25873333"; NickMod = "Ross PS Mod" },
    @{ Type = "Business Checking"; Nick = "Ross Corp"; Acct = "25874444"; NickMod = "Ross Corp Mod" }
)

for ($r = 0; $r -lt $achRows.Length; $r++) {
    $row = $r + 2
    $data = $achRows[$r]

    $ach.Cells.Item($row, 3).Value = "Required Fields ACH"
    $ach.Cells.Item($row, 4).Value = "Y"
    $ach.Cells.Item($row, 5).Value = "742"
    $ach.Cells.Item($row, 6).Value = "a_Access AutoNoCFtp"
    $ach.Cells.Item($row, 7).Value = "4249"
    $ach.Cells.Item($row, 8).Value = "a_Access AutoNoCFtpDemo"

    $ach.Cells.Item($row, 10).Value = "Delta Corp"
    $ach.Cells.Item($row, 11).Value = "Mr."
    $ach.Cells.Item($row, 12).Value = "Ross"
    $ach.Cells.Item($row, 13).Value = "KT"
    $ach.Cells.Item($row, 14).Value = "Evan"
    $ach.Cells.Item($row, 15).Value = "Sr."
    $ach.Cells.Item($row, 16).Value = "256987 Nolm Ct"
    $ach.Cells.Item($row, 17).Value = "Suite 678"
    $ach.Cells.Item($row, 18).Value = "United States"
    $ach.Cells.Item($row, 19).Value = "21054"
    $ach.Cells.Item($row, 20).Value = "1234"

    $ach.Cells.Item($row, 23).Value = "iahmed@govolution.com"
    $ach.Cells.Item($row, 24).Value = "240"
    $ach.Cells.Item($row, 25).Value = "628"
    $ach.Cells.Item($row, 26).Value = "0790"
    $ach.Cells.Item($row, 27).Value = "240"
    $ach.Cells.Item($row, 28).Value = "628"
    $ach.Cells.Item($row, 29).Value = "0791"
    $ach.Cells.Item($row, 30).Value = "410"
    $ach.Cells.Item($row, 31).Value = "628"
    $ach.Cells.Item($row, 32).Value = "0792"
    $ach.Cells.Item($row, 33).Value = "123"
    $ach.Cells.Item($row, 34).Value = "This Profile is for Add ACH"

    $ach.Cells.Item($row, 35).Value = $data.Type
    $ach.Cells.Item($row, 36).Value = "256072691"
    $ach.Cells.Item($row, 37).Value = $data.Nick
    $ach.Cells.Item($row, 38).Value = $data.Acct
    $ach.Cells.Item($row, 39).Value = $data.NickMod
}

$ach.Range("AE1").Select() | Out-Null
$ach.Range("AM1").Select() | Out-Null
$ach.Activate()

Write-Host "Added AddACHMRF and AddModifyDeleteACH sheets"
